$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell A2 with the slightly adjusted timestamp value
$ws.Range("A2").Value = 45804.53523658565

# Add new row 3 with new price data
$ws.Range("A3").Value = 45804.45374693675
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B3").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C3").Value = "1Kg"
$ws.Range("D3").Value = "12,88€"
